$d = $word.ActiveDocument

# --- Change 1: "VotingMachine" -> "PerfectCandidate" --------------------
# In Word, selecting "VotingMachine" and typing over it splits the
# enclosing run into three runs (before / new text / after), all sharing
# the original run formatting. Find.Execute repositions the range onto
# the matched text; setting .Text replaces it in place, and toggling a
# character-formatting property on that sub-range forces the run split
# while leaving formatting unchanged.
$hit = $d.Content
if ($hit.Find.Execute("VotingMachine")) {
    $start = $hit.Start
    $hit.Text = "PerfectCandidate"
    $newRange = $d.Range($start, $start + "PerfectCandidate".Length)
    $newRange.Font.Bold = $true
    $newRange.Font.Bold = $false
}

# --- Change 2: consolidate the "Once finished with voting..." runs -----
# The paragraph was previously split across four runs ("...refer to Step ",
# "3", ").", " Or, input...") with no formatting differences between them.
# Replacing the full matched text with itself collapses it back into a
# single run.
$full = "Once finished with voting, select the red exit button at the top left of the window. This will close the ballot and save the votes and winner of the election to another output file. Choosing which output file to write the information to has the same controls as finding the input file (refer to Step 3). Or, input an original file name to create a new file to be written to"
$d.Content.Find.Execute($full, $false, $false, $false, $false, $false, $true, 1, $false, $full, 2)
